$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.971.74"
$ws.Range("E2").Value = "  -5.07%  "

$ws.Range("D3").Value = "3.292.71"
$ws.Range("E3").Value = "  -6.20%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "'176.52"
$ws.Range("E5").Value = "  -11.77%  "

$ws.Range("D6").Value = "'523.78"
$ws.Range("E6").Value = "  -5.36%  "

$ws.Range("D7").Value = "'0.605"
$ws.Range("E7").Value = "  -0.45%  "

$ws.Range("D8").Value = "3.286.35"
$ws.Range("E8").Value = "  -6.16%  "

$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("D10").Value = "'0.607"
$ws.Range("E10").Value = "  -7.50%  "

$ws.Range("D11").Value = "'57.32"
$ws.Range("E11").Value = "  -7.80%  "

$ws.Range("E12").Value = "  -7.34%  "

$ws.Range("E13").Value = "  -4.82%  "

$ws.Range("D14").Value = "'9.07"
$ws.Range("E14").Value = "  -7.85%  "

$ws.Range("D15").Value = "3.816.11"
$ws.Range("E15").Value = "  -6.15%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.290.13"
$ws.Range("E16").Value = "  -6.03%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.117"
$ws.Range("E17").Value = "  -5.64%  "

$ws.Range("D18").Value = "63.813.05"
$ws.Range("E18").Value = "  -4.93%  "

$ws.Range("D19").Value = "'17.41"
$ws.Range("E19").Value = "  -5.81%  "

$ws.Range("D20").Value = "'11.07"
$ws.Range("E20").Value = "  -6.35%  "

$ws.Range("D21").Value = "'0.953"
$ws.Range("E21").Value = "  -7.28%  "

$ws.Range("D22").Value = "'373.39"
$ws.Range("E22").Value = "  -4.75%  "

$ws.Range("D23").Value = "'3.76"
$ws.Range("E23").Value = "  -6.12%  "

$ws.Range("E24").Value = "  -3.11%  "

$ws.Range("D25").Value = "'11.04"
$ws.Range("E25").Value = "  -10.28%  "

$ws.Range("D26").Value = "'3.86"
$ws.Range("E26").Value = "  -2.66%  "

$ws.Range("D27").Value = "'6.08"
$ws.Range("E27").Value = "  -1.36%  "

$ws.Range("D28").Value = "'2.66"
$ws.Range("E28").Value = "  -5.52%  "

$ws.Range("D29").Value = "'11.35"
$ws.Range("E29").Value = "  -7.31%  "

$ws.Range("D30").Value = "'8.33"
$ws.Range("E30").Value = "  -5.92%  "

$ws.Range("D31").Value = "'28.78"
$ws.Range("E31").Value = "  -7.15%  "

$ws.Range("D32").Value = "'636.89"
$ws.Range("E32").Value = "  -7.77%  "

$ws.Range("D33").Value = "'6.62"
$ws.Range("E33").Value = "  -6.01%  "

$ws.Range("D34").Value = "'11.22"
$ws.Range("E34").Value = "  -4.39%  "

$ws.Range("D35").Value = "'59.04"
$ws.Range("E35").Value = "  -7.51%  "

$ws.Range("E36").Value = "  -6.00%  "

$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("D38").Value = "'0.388"
$ws.Range("E38").Value = "  -2.19%  "

$ws.Range("D39").Value = "'36.49"
$ws.Range("E39").Value = "  -5.90%  "

$ws.Range("E40").Value = "  -0.09%  "

$ws.Range("E41").Value = "  +3.47%  "

$ws.Range("D42").Value = "2.922.71"
$ws.Range("E42").Value = "  -4.62%  "

$ws.Range("E43").Value = "  -5.21%  "

$ws.Range("D44").Value = "'2.46"
$ws.Range("E44").Value = "  -4.91%  "

$ws.Range("D45").Value = "'2.68"
$ws.Range("E45").Value = "  -10.48%  "

$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'2.65"
$ws.Range("E46").Value = "  -4.64%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0396"
$ws.Range("E47").Value = "  -1.75%  "

$ws.Range("D48").Value = "'3.03"

$ws.Range("D49").Value = "'2.77"
$ws.Range("E49").Value = "  +5.24%  "

$ws.Range("E50").Value = "  -1.50%  "

$ws.Range("D51").Value = "'135.37"
$ws.Range("E51").Value = "  -2.00%  "
